# The "data" sheet has a compact single-letter color-code column (E) whose
# header is "Color" and whose rows hold "y"/"g"/"r"/"b" (yellow/green/red/blue,
# matching the legend on 工作表2). Someone ran Find & Replace (Ctrl+H, "Replace
# All", searching the cell *contents* rather than whole cells) on the active
# "data" sheet to swap each single-letter code for its hex color, which also
# clipped the "r" out of the "Color" header turning it into "Colo#FF00FF".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# LookAt:=2 -> xlPart (substring match, not whole-cell) reproduces the
# "Color" -> "Colo#FF00FF" side effect seen in the diff.
$xlPart = 2
$xlByRows = 1

$ws.Cells.Replace("g", "#7FFF00", $xlPart, $xlByRows, $false, $false, $false, $false)
$ws.Cells.Replace("y", "#FFFF00", $xlPart, $xlByRows, $false, $false, $false, $false)
$ws.Cells.Replace("r", "#FF00FF", $xlPart, $xlByRows, $false, $false, $false, $false)
$ws.Cells.Replace("b", "#0000FF", $xlPart, $xlByRows, $false, $false, $false, $false)

# The user then scrolled down and clicked E10 before saving.
$ws.Range("E10").Select()
$excel.ActiveWindow.ScrollRow = 110
$excel.ActiveWindow.ScrollColumn = 1
